# Conserto do erro com o rotulo da coluna 2050 nas tabelas e retirada das
# linhas com total das tabelas.

$wb = $excel.ActiveWorkbook

# The year-header row (row 1) on five of the six tables had a stray leftover
# numeric value in E1 instead of the "2050" (or "2041-2050", for the period
# table) text label that belongs there alongside 2015/2030/2040 in B1:D1.
# Force the cell to stay text (like its neighbours) instead of letting a
# purely-numeric string be auto-coerced back into a number.
$sheetLabels = @(
    @{ Name = "Potencia Acumulada - SIN (MW)";  Label = "2050" },
    @{ Name = "Geracao Periodo Medio (MWMed)";  Label = "2050" },
    @{ Name = "Atendimento a Ponta(MW)";        Label = "2050" },
    @{ Name = "Potencia Incremental - SIN(MW)"; Label = "2041-2050" },
    @{ Name = "Emissoes Totais (MtCO2eq)";      Label = "2050" }
)

foreach ($entry in $sheetLabels) {
    $ws = $wb.Worksheets.Item($entry.Name)
    $cell = $ws.Range("E1")
    $cell.NumberFormat = "@"
    $cell.Value = $entry.Label
}

# Every table also carried a trailing "Total" row that should be dropped.
$totalRows = @(
    @{ Name = "Potencia Acumulada - SIN (MW)";  Row = 13 },
    @{ Name = "Geracao Periodo Medio (MWMed)";  Row = 13 },
    @{ Name = "Atendimento a Ponta(MW)";        Row = 13 },
    @{ Name = "Potencia Incremental - SIN(MW)"; Row = 13 },
    @{ Name = "Custo Total (bilhões de R$)";    Row = 4 }
)

foreach ($entry in $totalRows) {
    $ws = $wb.Worksheets.Item($entry.Name)
    $ws.Rows.Item($entry.Row).Delete()
}
